# =====================================================================
# Applies the "Add files via upload" commit:
#   * inserts a new blank "Sheet1" right after "디랙스"
#   * appends three new sheets: "고무블럭", "에버롤", "덤벨"
#   * adds a real hyperlink on 디랙스!C62 (pointing at its own image URL)
#   * adds real hyperlinks for every image-link cell on the new "고무블럭"
#     sheet
#   * leaves the final selection / active tab on "고무블럭" (matching the
#     author's last save state)
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New blank "Sheet1" right after "디랙스" (sheetId bumps to 6, so it
#    has to be created before the other new sheets).
# ---------------------------------------------------------------------
$drax = $wb.Worksheets.Item(1)
$blank = $wb.Worksheets.Add($null, $drax)
$blank.Name = "Sheet1"

# ---------------------------------------------------------------------
# 2) "고무블럭" (rubber block) sheet - appended at the very end, header +
#    17 product rows, each with a clickable hyperlink in column C.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$rubber = $wb.Worksheets.Add($null, $lastSheet)
$rubber.Name = "고무블럭"

$rubber.Columns.Item(1).ColumnWidth = 44.57142857142857
$rubber.Columns.Item(2).ColumnWidth = 29.857142857142858
$rubber.Columns.Item(3).ColumnWidth = 103.71428571428571

# header row - copy formatting from 디랙스!A1:C1 (style "1") so the new
# header cells land on the very same shared cellXf as every other sheet.
$rubber.Range("A1").Value = "제품"
$rubber.Range("B1").Value = "단가"
$rubber.Range("C1").Value = "사진"
$drax.Range("A1:C1").Copy()
$rubber.Range("A1:C1").PasteSpecial(-4122) | Out-Null

$rubberProducts = @(
    "일반 고무블럭 25T",
    "코팅 고무블럭 25T",
    "코팅 고무블럭 50T",
    "아레나 코팅 고무블럭 25T",
    "탑블랙 코팅 고무블럭 25T",
    "25T 마감재(실버) 일자 2400mm",
    "25T 마감재(실버) 모서리 300mm",
    "25T 마감재(실버) 역모서리 300mm",
    "50T 마감재(실버) 일자 2400mm",
    "50T 마감재(실버) 모서리 300mm",
    "50T 마감재(실버) 역모서리 300mm",
    "25T 마감재(골드) 일자 2400mm",
    "25T 마감재(골드) 모서리 300mm",
    "25T 마감재(골드) 역모서리 300mm",
    "25T 마감재(블랙) 일자 2400mm",
    "25T 마감재(블랙) 모서리 300mm",
    "25T 마감재(블랙) 역모서리 300mm"
)

$urlBase = "https://github.com/AsdDDsa182/AshGray/blob/main/images/RUBBERBLOCK/"

# use C2's url-link look (style "7", same as every DRAX image cell) as the
# formatting template for every C-column cell on this sheet.
$drax.Range("C2").Copy()

for ($i = 0; $i -lt $rubberProducts.Count; $i++) {
    $row = $i + 2
    $rubber.Cells.Item($row, 1).Value = $rubberProducts[$i]

    $cCell = $rubber.Cells.Item($row, 3)
    $cCell.PasteSpecial(-4122) | Out-Null
    $url = $urlBase + ($i + 1) + ".jpg?raw=true"
    $cCell.Value = $url
    $rubber.Hyperlinks.Add($cCell, $url) | Out-Null
}

# ---------------------------------------------------------------------
# 3) "에버롤" (Everoll) sheet - appended at the end, header + 1 row.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$everoll = $wb.Worksheets.Add($null, $lastSheet)
$everoll.Name = "에버롤"

$everoll.Columns.Item(1).ColumnWidth = 40
$everoll.Columns.Item(2).ColumnWidth = 46.42857142857143
$everoll.Columns.Item(3).ColumnWidth = 37.714285714285715

$everoll.Range("A1").Value = "제품"
$everoll.Range("B1").Value = "단가"
$everoll.Range("C1").Value = "사진"
$drax.Range("A1:C1").Copy()
$everoll.Range("A1:C1").PasteSpecial(-4122) | Out-Null

$everoll.Range("A2").Value = "에버롤 1x10m (1롤)"

# ---------------------------------------------------------------------
# 4) "덤벨" (Dumbbell) sheet - appended at the end, header + 1 row.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dumbbell = $wb.Worksheets.Add($null, $lastSheet)
$dumbbell.Name = "덤벨"

$dumbbell.Columns.Item(1).ColumnWidth = 43.285714285714285
$dumbbell.Columns.Item(2).ColumnWidth = 46.285714285714285
$dumbbell.Columns.Item(3).ColumnWidth = 41.714285714285715

$dumbbell.Range("A1").Value = "제품"
$dumbbell.Range("B1").Value = "단가"
$dumbbell.Range("C1").Value = "사진"
$drax.Range("A1:C1").Copy()
$dumbbell.Range("A1:C1").PasteSpecial(-4122) | Out-Null

$dumbbell.Range("A2").Value = "추가중"

# leftover selection state the author's workbook was saved with
$everoll.Range("B9").Select() | Out-Null
$dumbbell.Range("B9").Select() | Out-Null

# ---------------------------------------------------------------------
# 5) 디랙스!C62 keeps its text but now carries a real hyperlink to its
#    own image, and becomes the cell left selected on that sheet.
# ---------------------------------------------------------------------
$c62 = $drax.Range("C62")
$c62Url = $c62.Value
$drax.Hyperlinks.Add($c62, $c62Url) | Out-Null
$c62.Select() | Out-Null

# ---------------------------------------------------------------------
# 6) Final state: "고무블럭" is the active / tabSelected sheet, with
#    C12 selected (matches the author's last save).
# ---------------------------------------------------------------------
$rubber.Range("C12").Select() | Out-Null
